$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The recorded edit swaps the entire data content of row 2 and row 4
# (every populated column, A..AY) - the "Guldlockmossa" record that was
# on row 2 moves to row 4, and the "Mandelriska" record that was on row
# 4 moves to row 2. One column (L, "Kön") is present-but-blank only for
# whichever row holds the Guldlockmossa record, so its presence needs to
# travel with the swap too.

$colFirst = 1
$colLast = 51  # column AY

# Snapshot every value AND whether the cell exists at all (present-but-
# blank vs. truly absent) for row 2 and row 4 before writing anything.
# Value2 works reliably as a getter/setter here; plain Value does not
# (its getter misbehaves on this host, always returning a description
# string instead of the cell's contents).
$row2Vals = @{}
$row4Vals = @{}
$row2Exists = @{}
$row4Exists = @{}
for ($c = $colFirst; $c -le $colLast; $c++) {
    $v2 = $ws.Cells.Item(2, $c).Value2
    $v4 = $ws.Cells.Item(4, $c).Value2
    $row2Vals[$c] = $v2
    $row4Vals[$c] = $v4
    $row2Exists[$c] = ($null -ne $v2)
    $row4Exists[$c] = ($null -ne $v4)
}

for ($c = $colFirst; $c -le $colLast; $c++) {
    $v2 = $row2Vals[$c]
    $v4 = $row4Vals[$c]
    $e2 = $row2Exists[$c]
    $e4 = $row4Exists[$c]

    # Nothing to do if both rows already agree on this column (same
    # existence state and, when present, the same value) - e.g. the
    # Startdatum/Starttid/Slutdatum/Sluttid columns hold identical dates
    # in both rows. Skipping these also sidesteps Excel's automatic
    # text -> date coercion on Value2 writes for date-looking strings.
    if ($e2 -eq $e4 -and $v2 -eq $v4) {
        continue
    }

    $cell2 = $ws.Cells.Item(2, $c)
    $cell4 = $ws.Cells.Item(4, $c)

    $cell2.Value2 = $v4
    $cell4.Value2 = $v2

    # A blank-but-present cell (Value2 == "") needs a formatting touch to
    # force materialization - a bare Value2 = "" assignment leaves no
    # node behind (the cell reads back as completely absent) once the
    # value itself is empty.
    if ($e4 -and $v4 -eq "") {
        $cell2.Font.Bold = $cell2.Font.Bold
    }
    if ($e2 -and $v2 -eq "") {
        $cell4.Font.Bold = $cell4.Font.Bold
    }
}
